$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1), columns F through U ---
$ws.Range("F1").Value = "Loss function"
$ws.Range("G1").Value = "Avg. Time / Epoch"
$ws.Range("H1").Value = "Image dimension"
$ws.Range("I1").Value = "Loss"
$ws.Range("J1").Value = "Min. Loss"
$ws.Range("K1").Value = "Accuracy"
$ws.Range("L1").Value = "Dataset"
$ws.Range("M1").Value = "Device"
$ws.Range("N1").Value = "Convolutional layers"
$ws.Range("O1").Value = "Pools"
$ws.Range("P1").Value = "Created by"
$ws.Range("Q1").Value = "Total training time"
$ws.Range("R1").Value = "Gamma"
$ws.Range("S1").Value = "Weight decay"
$ws.Range("T1").Value = "Scheduler"
$ws.Range("U1").Value = "Min. LR"

# --- Update data rows 2-11 ---
# Row 2
$ws.Range("A2").Value = "2024-1-5 11:43:54"
$ws.Range("F2").Value = "CEL"
$ws.Range("G2").Value = 14.8
$ws.Range("H2").Value = 32
$ws.Range("I2").Value = 1.375
$ws.Range("J2").Value = 1.2159
$ws.Range("K2").Value = 53.4049
$ws.Range("L2").Value = "FER2013"
$ws.Range("M2").Value = "cuda:0"
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = "Stationær"
$ws.Range("Q2").Value = 296.1
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.005
$ws.Range("T2").Value = ""
$ws.Range("U2").Value = 0

# Row 3
$ws.Range("A3").Value = "2024-1-5 11:44:9"
$ws.Range("F3").Value = "CEL"
$ws.Range("G3").Value = 14.9
$ws.Range("H3").Value = 32
$ws.Range("I3").Value = 1.0891
$ws.Range("J3").Value = 1.0844
$ws.Range("K3").Value = 52.8127
$ws.Range("L3").Value = "FER2013"
$ws.Range("M3").Value = "cuda:0"
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = "Stationær"
$ws.Range("Q3").Value = 298.8
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.005
$ws.Range("T3").Value = ""
$ws.Range("U3").Value = 0

# Row 4
$ws.Range("A4").Value = "2024-1-5 11:44:11"
$ws.Range("F4").Value = "CEL"
$ws.Range("G4").Value = 14.9
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 1.4579
$ws.Range("J4").Value = 1.2095
$ws.Range("K4").Value = 53.3874
$ws.Range("L4").Value = "FER2013"
$ws.Range("M4").Value = "cuda:0"
$ws.Range("N4").Value = 4
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = "Stationær"
$ws.Range("Q4").Value = 297.4
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.005
$ws.Range("T4").Value = ""
$ws.Range("U4").Value = 0

# Row 5
$ws.Range("A5").Value = "2024-1-5 11:44:14"
$ws.Range("F5").Value = "CEL"
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 32
$ws.Range("I5").Value = 1.2336
$ws.Range("J5").Value = 1.1989
$ws.Range("K5").Value = 52.858
$ws.Range("L5").Value = "FER2013"
$ws.Range("M5").Value = "cuda:0"
$ws.Range("N5").Value = 4
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = "Stationær"
$ws.Range("Q5").Value = 299.3
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.005
$ws.Range("T5").Value = ""
$ws.Range("U5").Value = 0

# Row 6
$ws.Range("A6").Value = "2024-1-5 11:44:16"
$ws.Range("F6").Value = "CEL"
$ws.Range("G6").Value = 14.9
$ws.Range("H6").Value = 32
$ws.Range("I6").Value = 0.985
$ws.Range("J6").Value = 0.985
$ws.Range("K6").Value = 53.7183
$ws.Range("L6").Value = "FER2013"
$ws.Range("M6").Value = "cuda:0"
$ws.Range("N6").Value = 4
$ws.Range("O6").Value = 2
$ws.Range("P6").Value = "Stationær"
$ws.Range("Q6").Value = 298.9
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0.005
$ws.Range("T6").Value = ""
$ws.Range("U6").Value = 0

# Row 7
$ws.Range("A7").Value = "2024-1-5 11:44:18"
$ws.Range("F7").Value = "CEL"
$ws.Range("G7").Value = 15
$ws.Range("H7").Value = 32
$ws.Range("I7").Value = 1.1876
$ws.Range("J7").Value = 1.153
$ws.Range("K7").Value = 52.3111
$ws.Range("L7").Value = "FER2013"
$ws.Range("M7").Value = "cuda:0"
$ws.Range("N7").Value = 4
$ws.Range("O7").Value = 2
$ws.Range("P7").Value = "Stationær"
$ws.Range("Q7").Value = 299.4
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0.005
$ws.Range("T7").Value = ""
$ws.Range("U7").Value = 0

# Row 8
$ws.Range("A8").Value = "2024-1-5 11:44:19"
$ws.Range("F8").Value = "CEL"
$ws.Range("G8").Value = 14.9
$ws.Range("H8").Value = 32
$ws.Range("I8").Value = 1.2369
$ws.Range("J8").Value = 1.2369
$ws.Range("K8").Value = 52.5967
$ws.Range("L8").Value = "FER2013"
$ws.Range("M8").Value = "cuda:0"
$ws.Range("N8").Value = 4
$ws.Range("O8").Value = 2
$ws.Range("P8").Value = "Stationær"
$ws.Range("Q8").Value = 298.7
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0.005
$ws.Range("T8").Value = ""
$ws.Range("U8").Value = 0

# Row 9
$ws.Range("A9").Value = "2024-1-5 11:44:22"
$ws.Range("F9").Value = "CEL"
$ws.Range("G9").Value = 15
$ws.Range("H9").Value = 32
$ws.Range("I9").Value = 1.1981
$ws.Range("J9").Value = 0.9833
$ws.Range("K9").Value = 51.6389
$ws.Range("L9").Value = "FER2013"
$ws.Range("M9").Value = "cuda:0"
$ws.Range("N9").Value = 4
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = "Stationær"
$ws.Range("Q9").Value = 300.5
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0.005
$ws.Range("T9").Value = ""
$ws.Range("U9").Value = 0

# Row 10
$ws.Range("A10").Value = "2024-1-5 11:44:23"
$ws.Range("F10").Value = "CEL"
$ws.Range("G10").Value = 15
$ws.Range("H10").Value = 32
$ws.Range("I10").Value = 1.4878
$ws.Range("J10").Value = 1.0866
$ws.Range("K10").Value = 53.4954
$ws.Range("L10").Value = "FER2013"
$ws.Range("M10").Value = "cuda:0"
$ws.Range("N10").Value = 4
$ws.Range("O10").Value = 2
$ws.Range("P10").Value = "Stationær"
$ws.Range("Q10").Value = 299.2
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0.005
$ws.Range("T10").Value = ""
$ws.Range("U10").Value = 0

# Row 11
$ws.Range("A11").Value = "2024-1-5 11:44:24"
$ws.Range("F11").Value = "CEL"
$ws.Range("G11").Value = 15
$ws.Range("H11").Value = 32
$ws.Range("I11").Value = 1.491
$ws.Range("J11").Value = 1.2555
$ws.Range("K11").Value = 52.2415
$ws.Range("L11").Value = "FER2013"
$ws.Range("M11").Value = "cuda:0"
$ws.Range("N11").Value = 4
$ws.Range("O11").Value = 2
$ws.Range("P11").Value = "Stationær"
$ws.Range("Q11").Value = 299
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0.005
$ws.Range("T11").Value = "None"
$ws.Range("U11").Value = 0

